$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2295081967213115
$ws.Range("C2").Value = 0.5136612021857924
$ws.Range("J2").Value = 0.00819672131147541
$ws.Range("P2").Value = 0.185792349726776
$ws.Range("S2").Value = 0.06284153005464481
$ws.Range("B3").Value = 0.02040816326530612
$ws.Range("C3").Value = 0.03061224489795918
$ws.Range("J3").Value = 0.02551020408163265
$ws.Range("P3").Value = 0.7653061224489796
$ws.Range("S3").Value = 0.1581632653061225
$ws.Range("B6").Value = 0.04721030042918455
$ws.Range("D6").Value = 0.004291845493562232
$ws.Range("F6").Value = 0.04721030042918455
$ws.Range("J6").Value = 0.184549356223176
$ws.Range("O6").Value = 0.02575107296137339
$ws.Range("Q6").Value = 0.2274678111587983
$ws.Range("R6").Value = 0.09871244635193133
$ws.Range("S6").Value = 0.3648068669527897
$ws.Range("B7").Value = 0.1415094339622641
$ws.Range("D7").Value = 0.004716981132075472
$ws.Range("E7").Value = 0.004716981132075472
$ws.Range("F7").Value = 0.07547169811320754
$ws.Range("J7").Value = 0.1084905660377359
$ws.Range("O7").Value = 0.01415094339622642
$ws.Range("Q7").Value = 0.2075471698113208
$ws.Range("R7").Value = 0.09905660377358491
$ws.Range("S7").Value = 0.3443396226415094
$ws.Range("B8").Value = 0.09456740442655935
$ws.Range("D8").Value = 0.02012072434607646
$ws.Range("F8").Value = 0.04426559356136821
$ws.Range("J8").Value = 0.1146881287726358
$ws.Range("O8").Value = 0.01609657947686117
$ws.Range("Q8").Value = 0.1991951710261569
$ws.Range("R8").Value = 0.08249496981891348
$ws.Range("S8").Value = 0.4285714285714285
$ws.Range("B9").Value = 0.1333333333333333
$ws.Range("D9").Value = 0.009523809523809525
$ws.Range("F9").Value = 0.06666666666666667
$ws.Range("J9").Value = 0.1238095238095238
$ws.Range("O9").Value = 0.009523809523809525
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.08571428571428572
$ws.Range("S9").Value = 0.4047619047619048
$ws.Range("B10").Value = 0.114682249817385
$ws.Range("D10").Value = 0.02702702702702703
$ws.Range("F10").Value = 0.08327246165084003
$ws.Range("J10").Value = 0.1307523739956172
$ws.Range("O10").Value = 0.01241782322863404
$ws.Range("Q10").Value = 0.2045288531775018
$ws.Range("R10").Value = 0.08254200146092038
$ws.Range("S10").Value = 0.3447772096420745
$ws.Range("G11").Value = 0.14375
$ws.Range("J11").Value = 0.08125
$ws.Range("K11").Value = 0.19375
$ws.Range("L11").Value = 0.546875
$ws.Range("S11").Value = 0.034375
$ws.Range("G12").Value = 0.7471910112359551
$ws.Range("J12").Value = 0.1853932584269663
$ws.Range("K12").Value = 0.005617977528089887
$ws.Range("L12").Value = 0.02247191011235955
$ws.Range("S12").Value = 0.03932584269662921
$ws.Range("G13").Value = 0.6727272727272727
$ws.Range("J13").Value = 0.2727272727272727
$ws.Range("S13").Value = 0.05454545454545454
$ws.Range("F15").Value = 0.01435406698564593
$ws.Range("H15").Value = 0.1531100478468899
$ws.Range("I15").Value = 0.06698564593301436
$ws.Range("J15").Value = 0.3684210526315789
$ws.Range("K15").Value = 0.08133971291866028
$ws.Range("M15").Value = 0.01435406698564593
$ws.Range("O15").Value = 0.02870813397129187
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.008264462809917356
$ws.Range("H16").Value = 0.1694214876033058
$ws.Range("I16").Value = 0.1115702479338843
$ws.Range("J16").Value = 0.384297520661157
$ws.Range("K16").Value = 0.115702479338843
$ws.Range("M16").Value = 0.02066115702479339
$ws.Range("O16").Value = 0.05371900826446281
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.01377952755905512
$ws.Range("H17").Value = 0.1968503937007874
$ws.Range("I17").Value = 0.1003937007874016
$ws.Range("J17").Value = 0.4153543307086614
$ws.Range("K17").Value = 0.06692913385826772
$ws.Range("M17").Value = 0.02165354330708661
$ws.Range("O17").Value = 0.04527559055118111
$ws.Range("S17").Value = 0.1397637795275591
$ws.Range("F18").Value = 0.009216589861751152
$ws.Range("H18").Value = 0.2119815668202765
$ws.Range("I18").Value = 0.09677419354838709
$ws.Range("J18").Value = 0.3456221198156682
$ws.Range("K18").Value = 0.08294930875576037
$ws.Range("M18").Value = 0.03225806451612903
$ws.Range("O18").Value = 0.07373271889400922
$ws.Range("S18").Value = 0.1474654377880184
$ws.Range("F19").Value = 0.007315288953913679
$ws.Range("H19").Value = 0.2070226773957571
$ws.Range("I19").Value = 0.07388441843452817
$ws.Range("J19").Value = 0.3730797366495976
$ws.Range("K19").Value = 0.113386978785662
$ws.Range("M19").Value = 0.02267739575713241
$ws.Range("N19").Value = 0.000731528895391368
$ws.Range("O19").Value = 0.06583760058522312
$ws.Range("S19").Value = 0.1360643745427944
